$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7203
$ws.Range("C3").Value = 175452
$ws.Range("C4").Value = 165424
$ws.Range("C8").Value = 64.48999999999999
